$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.202518666666666
$ws.Cells.Item(2, 8).Value = 12.607556
$ws.Cells.Item(2, 9).Value = 0.08075097102331126
$ws.Cells.Item(2, 10).Value = 0.08075097102331129
$ws.Cells.Item(2, 13).Value = 16.57637
$ws.Cells.Item(2, 14).Value = 49.72911
$ws.Cells.Item(2, 15).Value = 0.1853914334114506
$ws.Cells.Item(2, 16).Value = 0.1853914334114506
$ws.Cells.Item(2, 17).Value = 69.66250435057333
$ws.Cells.Item(2, 18).Value = 626.9625391551599
$ws.Cells.Item(2, 19).Value = 0.01497053826737819
$ws.Cells.Item(2, 20).Value = 0.01497053826737819

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.202518666666666
$ws.Cells.Item(3, 8).Value = 12.607556
$ws.Cells.Item(3, 9).Value = 0.08075097102331126
$ws.Cells.Item(3, 10).Value = 0.08075097102331129
$ws.Cells.Item(3, 15).Value = 0.5978024790674488
$ws.Cells.Item(3, 16).Value = 0.5978024790674489
$ws.Cells.Item(3, 17).Value = 224.6296769624498
$ws.Cells.Item(3, 18).Value = 2021.667092662048
$ws.Cells.Item(3, 19).Value = 0.0482731306648392
$ws.Cells.Item(3, 20).Value = 0.04827313066483922

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.202518666666666
$ws.Cells.Item(4, 8).Value = 12.607556
$ws.Cells.Item(4, 9).Value = 0.08075097102331126
$ws.Cells.Item(4, 10).Value = 0.08075097102331129
$ws.Cells.Item(4, 15).Value = 0.2168060875211005
$ws.Cells.Item(4, 16).Value = 0.2168060875211005
$ws.Cells.Item(4, 17).Value = 81.46684416452977
$ws.Cells.Item(4, 18).Value = 733.2015974807679
$ws.Cells.Item(4, 19).Value = 0.01750730209109387
$ws.Cells.Item(4, 20).Value = 0.01750730209109388

$ws.Cells.Item(5, 9).Value = 0.7496282157262072
$ws.Cells.Item(5, 10).Value = 0.7496282157262073
$ws.Cells.Item(5, 13).Value = 16.57637
$ws.Cells.Item(5, 14).Value = 49.72911
$ws.Cells.Item(5, 15).Value = 0.1853914334114506
$ws.Cells.Item(5, 16).Value = 0.1853914334114506
$ws.Cells.Item(5, 17).Value = 646.69165184731
$ws.Cells.Item(5, 18).Value = 5820.22486662579
$ws.Cells.Item(5, 19).Value = 0.1389746494391497
$ws.Cells.Item(5, 20).Value = 0.1389746494391497

$ws.Cells.Item(6, 9).Value = 0.7496282157262072
$ws.Cells.Item(6, 10).Value = 0.7496282157262073
$ws.Cells.Item(6, 15).Value = 0.5978024790674488
$ws.Cells.Item(6, 16).Value = 0.5978024790674489
$ws.Cells.Item(6, 19).Value = 0.448129605740035
$ws.Cells.Item(6, 20).Value = 0.4481296057400351

$ws.Cells.Item(7, 9).Value = 0.7496282157262072
$ws.Cells.Item(7, 10).Value = 0.7496282157262073
$ws.Cells.Item(7, 15).Value = 0.2168060875211005
$ws.Cells.Item(7, 16).Value = 0.2168060875211005
$ws.Cells.Item(7, 19).Value = 0.1625239605470225
$ws.Cells.Item(7, 20).Value = 0.1625239605470225

$ws.Cells.Item(8, 7).Value = 8.827567333333333
$ws.Cells.Item(8, 9).Value = 0.1696208132504815
$ws.Cells.Item(8, 10).Value = 0.1696208132504815
$ws.Cells.Item(8, 13).Value = 16.57637
$ws.Cells.Item(8, 14).Value = 49.72911
$ws.Cells.Item(8, 15).Value = 0.1853914334114506
$ws.Cells.Item(8, 16).Value = 0.1853914334114506
$ws.Cells.Item(8, 17).Value = 146.3290223172467
$ws.Cells.Item(8, 18).Value = 1316.96120085522
$ws.Cells.Item(8, 19).Value = 0.03144624570492274
$ws.Cells.Item(8, 20).Value = 0.03144624570492275

$ws.Cells.Item(9, 7).Value = 8.827567333333333
$ws.Cells.Item(9, 9).Value = 0.1696208132504815
$ws.Cells.Item(9, 10).Value = 0.1696208132504815
$ws.Cells.Item(9, 15).Value = 0.5978024790674488
$ws.Cells.Item(9, 16).Value = 0.5978024790674489
$ws.Cells.Item(9, 17).Value = 471.8440905876462
$ws.Cells.Item(9, 18).Value = 4246.596815288816
$ws.Cells.Item(9, 19).Value = 0.1013997426625746
$ws.Cells.Item(9, 20).Value = 0.1013997426625746

$ws.Cells.Item(10, 7).Value = 8.827567333333333
$ws.Cells.Item(10, 9).Value = 0.1696208132504815
$ws.Cells.Item(10, 10).Value = 0.1696208132504815
$ws.Cells.Item(10, 15).Value = 0.2168060875211005
$ws.Cells.Item(10, 16).Value = 0.2168060875211005
$ws.Cells.Item(10, 19).Value = 0.03677482488298413
$ws.Cells.Item(10, 20).Value = 0.03677482488298414
